$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.184887070376633
$ws.Range("D2").Value = 8.907179426559852
$ws.Range("E2").Value = 13.43293486100976
$ws.Range("F2").Value = 34.97279954082592
$ws.Range("G2").Value = 3.626802777701948
$ws.Range("I2").Value = 24.171070806651
$ws.Range("J2").Value = 9.875322075795388
$ws.Range("N2").Value = 18.99769123883711
$ws.Range("O2").Value = 26.2486747980449

$ws.Range("C3").Value = 3.173039412892535
$ws.Range("D3").Value = 8.918618001937
$ws.Range("E3").Value = 13.41202101459894
$ws.Range("F3").Value = 34.56710278368384
$ws.Range("G3").Value = 3.630387562387526
$ws.Range("I3").Value = 23.91486072157148
$ws.Range("J3").Value = 9.878801804244505
$ws.Range("N3").Value = 18.40031508502701
$ws.Range("O3").Value = 25.97430490216302

$ws.Range("C4").Value = 3.165543185052695
$ws.Range("D4").Value = 8.92689574445712
$ws.Range("E4").Value = 13.40211214639553
$ws.Range("F4").Value = 34.3260860907202
$ws.Range("G4").Value = 3.632704804535492
$ws.Range("I4").Value = 23.76328479053294
$ws.Range("J4").Value = 9.882766463232553
$ws.Range("N4").Value = 18.02485520896362
$ws.Range("O4").Value = 25.81204298359441

$ws.Range("C5").Value = 3.162432947957773
$ws.Range("D5").Value = 8.930584414495634
$ws.Range("E5").Value = 13.39881351554923
$ws.Range("F5").Value = 34.23001325366301
$ws.Range("G5").Value = 3.63367842063957
$ws.Range("I5").Value = 23.703025702737
$ws.Range("J5").Value = 9.88484161047321
$ws.Range("N5").Value = 17.86990355188765
$ws.Range("O5").Value = 25.74754968645229

$ws.Range("C6").Value = 3.161913154073516
$ws.Range("D6").Value = 8.93121596504305
$ws.Range("E6").Value = 13.39831047878094
$ws.Range("F6").Value = 34.21419303463594
$ws.Range("G6").Value = 3.633841862984546
$ws.Range("I6").Value = 23.69311276154815
$ws.Range("J6").Value = 9.885213932416603
$ws.Range("N6").Value = 17.84406337566597
$ws.Range("O6").Value = 25.73694103359396

$ws.Range("C7").Value = 3.165501463507233
$ws.Range("D7").Value = 8.926944214058945
$ws.Range("E7").Value = 13.40206466430574
$ws.Range("F7").Value = 34.32478159745972
$ws.Range("G7").Value = 3.632717816198642
$ws.Range("I7").Value = 23.76246591904276
$ws.Range("J7").Value = 9.882792589188512
$ws.Range("N7").Value = 18.02277304767603
$ws.Range("O7").Value = 25.8111665158508

$ws.Range("C8").Value = 3.180848187150995
$ws.Range("D8").Value = 8.910863099260215
$ws.Range("E8").Value = 13.42511651026996
$ws.Range("F8").Value = 34.83130573294149
$ws.Range("G8").Value = 3.628014766264768
$ws.Range("I8").Value = 24.08157950035037
$ws.Range("J8").Value = 9.876142409083496
$ws.Range("N8").Value = 18.79364780656867
$ws.Range("O8").Value = 26.15282793226591

$ws.Range("C9").Value = 3.209171912478681
$ws.Range("D9").Value = 8.8892809528718
$ws.Range("E9").Value = 13.49347747762701
$ws.Range("F9").Value = 35.88364097495475
$ws.Range("G9").Value = 3.619708880248592
$ws.Range("I9").Value = 24.74981802204079
$ws.Range("J9").Value = 9.877608563381706
$ws.Range("N9").Value = 20.2273683202997
$ws.Range("O9").Value = 26.86879216756163

$ws.Range("C10").Value = 3.228886812376412
$ws.Range("D10").Value = 8.879490553274769
$ws.Range("E10").Value = 13.55763621742037
$ws.Range("F10").Value = 36.68574067888016
$ws.Range("G10").Value = 3.614158490423028
$ws.Range("I10").Value = 25.26237089453138
$ws.Range("J10").Value = 9.887524230839526
$ws.Range("N10").Value = 21.22223697909767
$ws.Range("O10").Value = 27.41830099875483

$ws.Range("C11").Value = 3.237614892093265
$ws.Range("D11").Value = 8.876352657301815
$ws.Range("E11").Value = 13.58979817134244
$ws.Range("F11").Value = 37.05538728041381
$ws.Range("G11").Value = 3.611751841278867
$ws.Range("I11").Value = 25.49928966524162
$ws.Range("J11").Value = 9.893949522941311
$ws.Range("N11").Value = 21.66018057919901
$ws.Range("O11").Value = 27.672386500647

$ws.Range("C12").Value = 3.240885143119419
$ws.Range("D12").Value = 8.875353429215879
$ws.Range("E12").Value = 13.60239967141692
$ws.Range("F12").Value = 37.19591936505106
$ws.Range("G12").Value = 3.610857396270401
$ws.Range("I12").Value = 25.58946456604469
$ws.Range("J12").Value = 9.896657164325664
$ws.Range("N12").Value = 21.82377585682186
$ws.Range("O12").Value = 27.76910816928867

$ws.Range("C13").Value = 3.2401823926377
$ws.Range("D13").Value = 8.875560227433061
$ws.Range("E13").Value = 13.59966701647175
$ws.Range("F13").Value = 37.16563065336736
$ws.Range("G13").Value = 3.611049280982884
$ws.Range("I13").Value = 25.57002466897496
$ws.Range("J13").Value = 9.896061830545055
$ws.Range("N13").Value = 21.78864458690801
$ws.Range("O13").Value = 27.74825633861114

$ws.Range("C14").Value = 3.237884638579147
$ws.Range("D14").Value = 8.876266663110885
$ws.Range("E14").Value = 13.59082647073783
$ws.Range("F14").Value = 37.0669385622324
$ws.Range("G14").Value = 3.61167791662476
$ws.Range("I14").Value = 25.50669967607238
$ws.Range("J14").Value = 9.894166785622023
$ws.Range("N14").Value = 21.67368539489659
$ws.Range("O14").Value = 27.68033422296782

$ws.Range("C15").Value = 3.236472645432399
$ws.Range("D15").Value = 8.876723985602524
$ws.Range("E15").Value = 13.58546622723222
$ws.Range("F15").Value = 37.00655521686857
$ws.Range("G15").Value = 3.612065172022333
$ws.Range("I15").Value = 25.4679686444845
$ws.Range("J15").Value = 9.893041739964682
$ws.Range("N15").Value = 21.60297336126124
$ws.Range("O15").Value = 27.63879317500566

$ws.Range("C16").Value = 3.228311548051427
$ws.Range("D16").Value = 8.879722084771453
$ws.Range("E16").Value = 13.55559378609648
$ws.Range("F16").Value = 36.66166815196529
$ws.Range("G16").Value = 3.614318141252793
$ws.Range("I16").Value = 25.24695633307867
$ws.Range("J16").Value = 9.887142786062523
$ws.Range("N16").Value = 21.19330956972086
$ws.Range("O16").Value = 27.4017712756022

$ws.Range("C17").Value = 3.223243194034749
$ws.Range("D17").Value = 8.88189821450387
$ws.Range("E17").Value = 13.53802649215905
$ws.Range("F17").Value = 36.4512189864138
$ws.Range("G17").Value = 3.615730477979618
$ws.Range("I17").Value = 25.11227666083764
$ws.Range("J17").Value = 9.884013840549242
$ws.Range("N17").Value = 20.93814219015166
$ws.Range("O17").Value = 27.25735771348802

$ws.Range("C18").Value = 3.220305409422927
$ws.Range("D18").Value = 8.883273731021808
$ws.Range("E18").Value = 13.52820271101434
$ws.Range("F18").Value = 36.33063071707114
$ws.Range("G18").Value = 3.616553953218962
$ws.Range("I18").Value = 25.03517087522813
$ws.Range("J18").Value = 9.882394425401978
$ws.Range("N18").Value = 20.79000725568362
$ws.Range("O18").Value = 27.17468699768372

$ws.Range("C19").Value = 3.21930685413352
$ws.Range("D19").Value = 8.883760735918965
$ws.Range("E19").Value = 13.52492487384914
$ws.Range("F19").Value = 36.28988400255999
$ws.Range("G19").Value = 3.616834683654054
$ws.Range("I19").Value = 25.00912821015849
$ws.Range("J19").Value = 9.881877102131813
$ws.Range("N19").Value = 20.73962067985785
$ws.Range("O19").Value = 27.14676608934397

$ws.Range("C20").Value = 3.223785070506256
$ws.Range("D20").Value = 8.881653743502477
$ws.Range("E20").Value = 13.53986757158044
$ws.Range("F20").Value = 36.4735753601234
$ws.Range("G20").Value = 3.615578980493327
$ws.Range("I20").Value = 25.12657703735449
$ws.Range("J20").Value = 9.884328269894427
$ws.Range("N20").Value = 20.96544799484619
$ws.Range("O20").Value = 27.27269082254472

$ws.Range("C21").Value = 3.238560493944457
$ws.Range("D21").Value = 8.876054037480063
$ws.Range("E21").Value = 13.59341173398027
$ws.Range("F21").Value = 37.09591279659214
$ws.Range("G21").Value = 3.611492813237964
$ws.Range("I21").Value = 25.52528795337747
$ws.Range("J21").Value = 9.894715963076063
$ws.Range("N21").Value = 21.70751365554066
$ws.Range("O21").Value = 27.70027157568724

$ws.Range("C22").Value = 3.248013611458014
$ws.Range("D22").Value = 8.873495965859222
$ws.Range("E22").Value = 13.63086546988129
$ws.Range("F22").Value = 37.50582122806112
$ws.Range("G22").Value = 3.608920729752568
$ws.Range("I22").Value = 25.78850597668649
$ws.Range("J22").Value = 9.903104544899309
$ws.Range("N22").Value = 22.17935961385674
$ws.Range("O22").Value = 27.982624045944

$ws.Range("C23").Value = 3.242987029473977
$ws.Range("D23").Value = 8.874760529569944
$ws.Range("E23").Value = 13.61065260601627
$ws.Range("F23").Value = 37.28679771695708
$ws.Range("G23").Value = 3.610284523830297
$ws.Range("I23").Value = 25.64780700570665
$ws.Range("J23").Value = 9.898481347597302
$ws.Range("N23").Value = 21.92877110911181
$ws.Range("O23").Value = 27.83168995744535

$ws.Range("C24").Value = 3.223540162627626
$ws.Range("D24").Value = 8.881763881248403
$ws.Range("E24").Value = 13.53903435982324
$ws.Range("F24").Value = 36.46346678210779
$ws.Range("G24").Value = 3.615647436663764
$ws.Range("I24").Value = 25.12011082369111
$ws.Range("J24").Value = 9.884185557348149
$ws.Range("N24").Value = 20.95310750188673
$ws.Range("O24").Value = 27.26575761288651

$ws.Range("C25").Value = 3.201700739316915
$ws.Range("D25").Value = 8.89405394664173
$ws.Range("E25").Value = 13.47252147058638
$ws.Range("F25").Value = 35.59335385997608
$ws.Range("G25").Value = 3.621858415089816
$ws.Range("I25").Value = 24.56493140925214
$ws.Range("J25").Value = 9.87565867579619
$ws.Range("N25").Value = 19.84905939529497
$ws.Range("O25").Value = 26.670643588452

